# reDIP-RIOT-BOM.xlsx: swap 74AHC1G14 -> 74AHCT1G14 (U5) for better signal
# level margins, and refresh the KiCost-sourced pricing/stock/date data that
# comes along with re-running the part lookup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Part-number text (shared strings used by cells B19, D19, G19, P19, V19)
# ---------------------------------------------------------------------
$ws.Range("B19").Value = "74AHCT1G14"
$ws.Range("D19").Value = "Replacement parts: 74AHCT1G14 SOT-353/SC-70-5"
$ws.Range("G19").Value = "74AHCT1G14SE-7"
$ws.Range("P19").Value = "74AHCT1G14SE-7DICT-ND"
$ws.Range("V19").Value = "621-74AHCT1G14SE-7"

# ---------------------------------------------------------------------
# 2. Distributor #1 (Digi-Key columns K-O) pricing/stock
# ---------------------------------------------------------------------
$ws.Range("K19").Value = 74194
$ws.Range("M19").Formula = '=IFERROR(IF(OR(L19>=N19,H19>=N19),LOOKUP(IF(L19="",H19,L19),{0,1,10,25,100,250,500,1000},{0.0,0.27,0.223,0.1864,0.118,0.09112,0.07764,0.0528}),"MOQ="&N19),"")'

# ---------------------------------------------------------------------
# 3. Distributor #2 (Mouser columns Q-U) pricing/stock
# ---------------------------------------------------------------------
$ws.Range("Q19").Value = 1287
$ws.Range("S19").Formula = '=IFERROR(IF(OR(R19>=T19,H19>=T19),LOOKUP(IF(R19="",H19,R19),{0,1,10,100,1000,3000,9000,24000,45000},{0.0,0.269,0.184,0.077,0.052,0.041,0.034,0.032,0.03}),"MOQ="&T19),"")'

# ---------------------------------------------------------------------
# 4. KiCost run timestamps (B3 / B4)
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "fr. 06. okt. 2023 kl. 12.39 +0200"
$ws.Range("B4").Value = "2023-10-06 12:40:02"

# ---------------------------------------------------------------------
# 5. Cell comments (hover tooltips) for row 19
# ---------------------------------------------------------------------
$m19Comment = 'Qty/Price Breaks (USD):
  Qty  -  Unit$  -  Ext$
================
     1   $0.27      $0.27
    10   $0.22      $2.23
    25   $0.19      $4.66
   100   $0.12     $11.80
   250   $0.09     $22.78
   500   $0.08     $38.82
  1000   $0.05     $52.80'
$null = $ws.Range("M19").Comment.Text($m19Comment)

$p19Comment = 'Desc: IC INVERTER 1CH 1-INP SOT353
Footprint: 5-TSSOP, SC-70-5, SOT-353'
$null = $ws.Range("P19").Comment.Text($p19Comment)

$null = $ws.Range("Q19").Comment.Text("1287 In Stock")

$s19Comment = 'Qty/Price Breaks (USD):
  Qty  -  Unit$  -  Ext$
================
     1   $0.27      $0.27
    10   $0.18      $1.84
   100   $0.08      $7.70
  1000   $0.05     $52.00
  3000   $0.04    $123.00
  9000   $0.03    $306.00
 24000   $0.03    $768.00
 45000   $0.03  $1,350.00'
$null = $ws.Range("S19").Comment.Text($s19Comment)

$null = $ws.Range("V19").Comment.Text("Desc: Inverters AHC TTL Compatible LOGIC")
